# Daily attendance processing - 2026-01-04 23:34:06
# Swap the order of "dnasr281@gmail.com" and "System" in the
# "Recorded By" column (G) wherever both appear together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "dnasr281@gmail.com, System"
$newText = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()
    if ($val -eq $oldText) {
        $cell.Value = $newText
    }
}
